$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells L1:N1, matching the style of the existing headers (K1) ---
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- Rescale columns E (particip) and F (taxa_sucesso) from fraction to 0-100 scale ---
$particip = @(91.08614232209737, 8.913857677902621, 92.8474114441417, 7.152588555858311, 94.88304093567251, 5.116959064327485)
$taxa = @(60.9375, 74.78991596638656, 93.91049156272928, 98.09523809523809, 20.80123266563945, 48.57142857142857)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = $particip[$i]
    $ws.Cells.Item($r, 6).Value = $taxa[$i]
}

# --- New data columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes) ---
$apoioMedio = @(92.19084927254757, 82.97112949625762, 90.90970086312072, 81.91974458046401, 19.81586961574516, 14.52164198228234)
$contribuicoes = @(238184, 25369, 186820, 16826, 2101, 107)
$mediaContribuicoes = @(321.4358974358975, 285.0449438202247, 145.953125, 163.3592233009709, 15.56296296296296, 6.294117647058823)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 12).Value = $apoioMedio[$i]
    $ws.Cells.Item($r, 13).Value = $contribuicoes[$i]
    $ws.Cells.Item($r, 14).Value = $mediaContribuicoes[$i]
}
